# Re-apply the table style used by the three summary tables (slides 14-16)
# from the default "Table_0" style ({CA9A5765-2B38-4931-AF24-86F3967F88BB})
# to the PowerPoint built-in table style {CDD7A4D9-40D5-4805-A191-EDB4B8677A5A}.
#
# Table styles can only be changed via Table.ApplyStyle(id) - the Style /
# TableStyleId properties are read-only on the Table object.

$p = $ppt.ActivePresentation

$oldStyleId = "{CA9A5765-2B38-4931-AF24-86F3967F88BB}"
$newStyleId = "{CDD7A4D9-40D5-4805-A191-EDB4B8677A5A}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.TableStyleId -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
